$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.651.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "'2.124.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "'353.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.48%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "'0.5282"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").Value = "'0.4536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'54.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.09108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'1.183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'2.115.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'8.117"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "'102.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.15%  "
$ws.Range("D17").Value = "'0.00001178"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'0.06719"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'19.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'6.350"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'30.727.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").Value = "'2.398"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").Value = "'2.371.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'2.577"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").Value = "'165.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'136.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("D31").Value = "'1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'0.1081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'1.661"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "'6.395"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "'4.025"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D38").Value = "'0.02663"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "'0.06897"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "'12.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "'0.6936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'1.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").Value = "'14.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.56%  "
$ws.Range("D47").Value = "'3.769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("E48").Value = "  +8.75%  "
$ws.Range("D49").Value = "'1.259"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

# Row 36/37 swap: FraxShare <-> InternetComputer(DFINITY)
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'6.025"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.82%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'10.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.23%  "

# Row 45/46 swap: Decentraland <-> NEARProtocol
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.340"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6485"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "

# Row 50/51 swap: Cronos <-> Aave
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'83.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.07319"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.50%  "
